# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
# D-column values that look like plain decimal numbers must be forced to
# stay text (they are formatted/truncated strings, not real numbers), so
# those are written with a leading apostrophe (Excel's "treat as text"
# quote-prefix). D-column values containing two dots (e.g. "26.839.62")
# are never auto-converted by Excel, so they are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.839.62";  DText = $false; E = "  -1.30%  " }
    @{ Row = 3;  D = "1.874.13";   DText = $false; E = "  -1.55%  " }
    @{ Row = 4;  E = "  -0.17%  " }
    @{ Row = 5;  D = "301.78";     DText = $true;  E = "  -1.86%  " }
    @{ Row = 6;  E = "  -0.10%  " }
    @{ Row = 7;  D = "0.5367";     DText = $true;  E = "  +2.29%  " }
    @{ Row = 8;  D = "0.3760";     DText = $true;  E = "  -1.54%  " }
    @{ Row = 9;  D = "0.07190";    DText = $true }
    @{ Row = 10; D = "21.59";      DText = $true;  E = "  -0.34%  " }
    @{ Row = 12; D = "0.08165";    DText = $true;  E = "  +0.70%  " }
    @{ Row = 13; D = "1.872.53";   DText = $false; E = "  +2.92%  " }
    @{ Row = 14; D = "93.51";      DText = $true;  E = "  -2.23%  " }
    @{ Row = 15; D = "5.268";      DText = $true;  E = "  -1.70%  " }
    @{ Row = 16; D = "1.001";      DText = $true;  E = "  -0.19%  " }
    @{ Row = 17; D = "14.75";      DText = $true;  E = "  +0.12%  " }
    @{ Row = 19; E = "  -0.01%  " }
    @{ Row = 20; D = "26.887.63";  DText = $false; E = "  -1.28%  " }
    @{ Row = 21; D = "4.983";      DText = $true;  E = "  -2.58%  " }
    @{ Row = 22; D = "10.68";      DText = $true;  E = "  -1.19%  " }
    @{ Row = 23; D = "6.396";      DText = $true;  E = "  -1.15%  " }
    @{ Row = 24; D = "147.06";     DText = $true;  E = "  -1.57%  " }
    @{ Row = 25; D = "2.257";      DText = $true;  E = "  -3.32%  " }
    @{ Row = 26; D = "1.736";      DText = $true;  E = "  -0.44%  " }
    @{ Row = 27; E = "  -1.20%  " }
    @{ Row = 28; D = "114.02";     DText = $true;  E = "  -1.99%  " }
    @{ Row = 29; D = "4.727";      DText = $true;  E = "  -2.26%  " }
    @{ Row = 30; D = "4.604";      DText = $true;  E = "  -5.91%  " }
    @{ Row = 31; D = "0.09164";    DText = $true }
    @{ Row = 32; D = "0.8052";     DText = $true;  E = "  +0.96%  " }
    @{ Row = 33; D = "0.04978";    DText = $true;  E = "  -2.00%  " }
    @{ Row = 34; D = "1.173";      DText = $true;  E = "  -4.58%  " }
    @{ Row = 35; D = "2.984";      DText = $true;  E = "  +0.10%  " }
    @{ Row = 36; D = "0.6023";     DText = $true;  E = "  +5.00%  " }
    @{ Row = 37; D = "3.203";      DText = $true;  E = "  -5.05%  " }
    @{ Row = 38; D = "2.561";      DText = $true;  E = "  -4.29%  " }
    @{ Row = 39; D = "0.01951";    DText = $true;  E = "  -2.30%  " }
    @{ Row = 40; D = "1.072";      DText = $true;  E = "  -1.21%  " }
    @{ Row = 41; D = "6.570";      DText = $true;  E = "  -0.46%  " }
    @{ Row = 42; D = "8.855";      DText = $true;  E = "  -1.89%  " }
    @{ Row = 43; D = "0.5145";     DText = $true;  E = "  +5.13%  " }
    @{ Row = 44; D = "115.06";     DText = $true;  E = "  -1.48%  " }
    @{ Row = 45; E = "  -1.65%  " }
    @{ Row = 46; E = "  -0.14%  " }
    @{ Row = 47; D = "9.923";      DText = $true;  E = "  -2.61%  " }
    @{ Row = 48; E = "  -0.08%  " }
    @{ Row = 49; D = "37.56";      DText = $true;  E = "  -2.68%  " }
    @{ Row = 50; D = "0.06026";    DText = $true;  E = "  +1.10%  " }
    @{ Row = 51; E = "  -3.06%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("D")) {
        $dCell = $ws.Range("D$row")
        if ($u.DText) {
            # Leading apostrophe forces Excel to store the value as text
            # instead of silently re-parsing it as a number (which would
            # both change the cell type and drop formatting like trailing
            # zeros, e.g. "0.3760" -> 0.376).
            $dCell.Value = "'" + $u.D
        } else {
            $dCell.Value = $u.D
        }
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
